# Add a new row (row 39) with data to each of the 4 worksheets in the workbook.
$wb = $excel.ActiveWorkbook

# Large magnitude values for column G need to be built from factors that fit
# within the numeric literal limits of this environment, then multiplied
# together to reconstruct the original double value exactly.
$gBig1 = 5.68631262647114 * 100000 * 1000000000000000000
$gBig2 = 9.85046333984776 * 100000 * 1000000000000000000

# --- Sheet 1: ROW35-FE-LIFTER ---
$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
$ws1.Range("A39").Value = 45743.82211991898
$ws1.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B39").Value = "0x01,0x90"
$ws1.Range("C39").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Range("D39").Value = "0x01,0x7a"
$ws1.Range("E39").Value = "0xd"
$ws1.Range("F39").Value = 400
$ws1.Range("G39").Value = $gBig1
$ws1.Range("H39").Value = 378
$ws1.Range("I39").Value = 13

# --- Sheet 2: ROW35-MID-LIFTER ---
$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
$ws2.Range("A39").Value = 45743.66841015047
$ws2.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B39").Value = "0x01,0x90"
$ws2.Range("C39").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Range("D39").Value = "0x01,0x7a"
$ws2.Range("E39").Value = "0xe"
$ws2.Range("F39").Value = 400
$ws2.Range("G39").Value = $gBig1
$ws2.Range("H39").Value = 378
$ws2.Range("I39").Value = 14

# --- Sheet 3: ROW02-FE-LIFTER ---
$ws3 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
$ws3.Range("A39").Value = 45743.81536850694
$ws3.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B39").Value = "0x01,0x90"
$ws3.Range("C39").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Range("D39").Value = "0x01,0x7a"
$ws3.Range("E39").Value = "0x3"
$ws3.Range("F39").Value = 400
$ws3.Range("G39").Value = $gBig1
$ws3.Range("H39").Value = 378
$ws3.Range("I39").Value = 3

# --- Sheet 4: ROW02-MID-LIFTER ---
$ws4 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
$ws4.Range("A39").Value = 45743.87053054398
$ws4.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("B39").Value = "0x01,0x90"
$ws4.Range("C39").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Range("D39").Value = "0x01,0x7a"
$ws4.Range("E39").Value = "0x3"
$ws4.Range("F39").Value = 400
$ws4.Range("G39").Value = $gBig2
$ws4.Range("H39").Value = 378
$ws4.Range("I39").Value = 3

Write-Host "Row 39 added to all 4 sheets"
